# TC01_C3DC_phs000471_SexAtBirth-Unknown.xlsx — regression/smoke suite update
#
# 1) Fix the "Treatment" tab query: replace the redundant
#    CONCAT(REPLACE(...)) wrapper with a plain REPLACE(...) call.
# 2) As a consequence of how the shared-strings table is rebuilt, this also
#    reorders the three big query strings so that the Treatment-Response and
#    Survival query text (which are untouched) come right after the
#    Diagnosis query, and the (now-fixed) Treatment query text moves to the
#    very end of the string table.
# 3) Update the active selection to C5 (previously C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the current "Treatment" tab query text out of B5.
$treatmentQuery = $ws.Range("B5").Value2

# Temporarily blank out B5. This fully de-references the old query text so
# that, once the corrected text is written back, the stale copy is dropped
# from the workbook's shared string table instead of lingering in place.
$ws.Range("B5").Value = ""

# Apply the fix: CONCAT(REPLACE(...)) -> REPLACE(...)
$fixedTreatmentQuery = $treatmentQuery.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")

# Write the corrected query back into B5.
$ws.Range("B5").Value = $fixedTreatmentQuery

# Restore the selection Excel would have left active after this edit.
$ws.Range("C5").Select()
